$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'tight calf muscles',
    'tight calves',
    'tight capri',
    'tight capris for women',
    'tight clothes',
    'tight clothes for women',
    'tight compression shirt',
    'tight compression underwear',
    'tight eyes',
    'tight fit',
    'tight fit leggings women',
    'tight fitted pants for women',
    'tight hamstrings',
    'tight high stockings',
    'tight hip flexors',
    'tight hold',
    'tight hoops',
    'tight jacket',
    'tight jackets for women',
    'tight leggings',
    'tight leggings see through',
    'tight leggings women',
    'tight master',
    'tight muscle',
    'tight muscle therapy',
    'tight muscles',
    'tight overalls',
    'tight overalls for women',
    'tight pant',
    'tight panties',
    'tight pants',
    'tight pants for men',
    'tight pants for women',
    'tight pants women',
    'tight running shorts women',
    'tight shiny pants',
    'tight short',
    'tight shorts',
    'tight shorts for women',
    'tight shorts women',
    'tight ski pants women',
    'tight skin',
    'tight slimmer for men',
    'tight slimmer for women',
    'tight small',
    'tight spandex shorts',
    'tight sport women puma',
    'tight spot',
    'tight spring',
    'tight stick',
    'tight sweats for women',
    'tight teen',
    'tight tennis skirt',
    'tight top',
    'tight tops',
    'tight tummy',
    'tight tummy for her',
    'tight up skirt',
    'tight weave cloth',
    'tight white shorts',
    'tight women',
    'tight womens ski pants',
    'tight workout tops for women',
    'tight yoga shorts',
    'tighten loose skin',
    'tighten tummy skin',
    'tights 3x',
    'tights 5x',
    'tights and leggings',
    'tights basketball',
    'tights black',
    'tights brown',
    'tights compression',
    'tights extra long',
    'tights floral',
    'tights for',
    'tights for cold',
    'tights for cold weather women',
    'tights for running',
    'tights for running womens',
    'tights for soccer',
    'tights for teens',
    'tights for track',
    'tights for winter for women',
    'tights for women',
    'tights for women athletic',
    'tights for women black',
    'tights for women capri',
    'tights for women gray',
    'tights for women grey',
    'tights for women gym',
    'tights for women high waist',
    'tights for women knee high',
    'tights for women knee length',
    'tights for women pack',
    'tights for women plus size',
    'tights for women with designs',
    'tights for women with pockets',
    'tights for women workout',
    'tights for women xl'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
